# Huge bugfix that fixed the problem where we could not read the same file multiple times

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("task list")

# --- Update the book/window view: minimize the workbook window and move it ---
$win = $excel.ActiveWindow
$win.WindowState = -4128   # xlMinimized
$win.Left = 6400
$win.Top = 6060

# --- Update the status of the first four line items (rows 3-5 = Completed, row 6 = In Progress) ---
$ws.Range("E3").Value = "Completed"
$ws.Range("E4").Value = "Completed"
$ws.Range("E5").Value = "Completed"
$ws.Range("E6").Value = "In Progress"

# --- Move the active selection to F5 ---
$ws.Range("F5").Select()
